$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per diff
$ws.Range("B29").Value = "J"
$ws.Range("B59").Value = "Y"

# Update selected cell in sheet view
$ws.Range("G3").Select()
